{"js": "// Adds three \"ALTER TABLE ... ADD FOREIGN KEY ...\" statements to the\n// (currently empty) paragraphs that follow the existing\n// \"ALTER TABLE DEPARTMENTS...\", \"ALTER TABLE EMPLOYEES...\" and\n// \"ALTER TABLE JOB_HISTORY...\" paragraphs, reproducing:\n//   alter table DEPARTMENTS ADD FOREIGN KEY(location_id) REFERENCES LOCATIONS(location_id);\n//   alter table JOB_HISTORY ADD FOREIGN KEY(employee_id) REFERENCES EMPLOYEES(employee_id);\n//   alter table JOB_HISTORY ADD FOREIGN KEY(job_id) REFERENCES JOBS(job_id);\n// Each new run carries <w:lang w:val=\"en-US\"/> and the dictionary-flagged\n// identifiers (location_id / employee_id / job_id) are wrapped in\n// w:proofErr spellStart/spellEnd markers, matching Word's own \"as you\n// type\" behaviour for an unrecognised word.\n\n// Build a <pkg:package> OOXML fragment for one paragraph's worth of runs.\n// `segments` is an array of [text, isFlaggedWord] pairs.\nfunction buildRunsOoxml(segments) {\n  const nsDecl =\n    'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n  let runs = \"\";\n  for (const [text, flagged] of segments) {\n    const escaped = text\n      .replace(/&/g, \"&amp;\")\n      .replace(/</g, \"&lt;\")\n      .replace(/>/g, \"&gt;\");\n    const preserve = /^\\s|\\s$/.test(text) ? ' xml:space=\"preserve\"' : \"\";\n    const run =\n      \"<w:r><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr>\" +\n      \"<w:t\" + preserve + \">\" + escaped + \"</w:t></w:r>\";\n    if (flagged) {\n      runs += '<w:proofErr w:type=\"spellStart\"/>' + run + '<w:proofErr w:type=\"spellEnd\"/>';\n    } else {\n      runs += run;\n    }\n  }\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    \"<w:document \" + nsDecl + \"><w:body><w:p>\" + runs + \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Locate the (empty) paragraph that immediately follows the paragraph\n// whose text starts with `precedingText`, and insert the OOXML for\n// `segments` into it.\nasync function fillParagraphAfter(context, precedingText, segments) {\n  // matchCase:true \u2014 the sentinel text is the pre-existing, ALL-CAPS\n  // \"ALTER TABLE ...\" statement; the newly inserted statements below all\n  // start with lower-case \"alter table\", so a case-sensitive search keeps\n  // later calls from matching text this function itself just inserted.\n  const results = context.document.body.search(precedingText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const anchorParagraph = results.items[0].paragraphs.getFirst();\n  const target = anchorParagraph.getNext();\n  target.insertOoxml(buildRunsOoxml(segments), Word.InsertLocation.end);\n  await context.sync();\n}\n\nawait fillParagraphAfter(\n  context,\n  \"ALTER TABLE DEPARTMENTS ADD FOREIGN KEY\",\n  [\n    [\"alter table DEPARTMENTS ADD FOREIGN KEY(\", false],\n    [\"location_id\", true],\n    [\") REFERENCES LOCATIONS(\", false],\n    [\"location_id\", true],\n    [\");\", false],\n  ]\n);\n\nawait fillParagraphAfter(\n  context,\n  \"ALTER TABLE EMPLOYEES ADD FOREIGN KEY\",\n  [\n    [\"alter table JOB_HISTORY ADD FOREIGN KEY(\", false],\n    [\"employee_id\", true],\n    [\") REFERENCES EMPLOYEES(\", false],\n    [\"employee_id\", true],\n    [\");\", false],\n  ]\n);\n\nawait fillParagraphAfter(\n  context,\n  \"ALTER TABLE JOB_HISTORY ADD FOREIGN KEY\",\n  [\n    [\"alter table JOB_HISTORY ADD FOREIGN KEY(\", false],\n    [\"job_id\", true],\n    [\") REFERENCES JOBS(\", false],\n    [\"job_id\", true],\n    [\");\", false],\n  ]\n);\n", "ps1": "# Adds three \"ALTER TABLE ... ADD FOREIGN KEY ...\" statements to the\n# (currently empty) paragraphs that follow the existing\n# \"ALTER TABLE DEPARTMENTS...\", \"ALTER TABLE EMPLOYEES...\" and\n# \"ALTER TABLE JOB_HISTORY...\" paragraphs, reproducing:\n#   alter table DEPARTMENTS ADD FOREIGN KEY(location_id) REFERENCES LOCATIONS(location_id);\n#   alter table JOB_HISTORY ADD FOREIGN KEY(employee_id) REFERENCES EMPLOYEES(employee_id);\n#   alter table JOB_HISTORY ADD FOREIGN KEY(job_id) REFERENCES JOBS(job_id);\n# Each new run carries <w:lang w:val=\"en-US\"/> and the dictionary-flagged\n# identifiers (location_id / employee_id / job_id) are wrapped in\n# w:proofErr spellStart/spellEnd markers, matching Word's own \"as you\n# type\" behaviour for an unrecognised word.\n\n$d = $word.ActiveDocument\n\n# Build the <w:r>...</w:r> (optionally w:proofErr-wrapped) run sequence for\n# one paragraph's worth of text. $Segments is an array of 2-element arrays:\n# @(text, isFlaggedWord).\nfunction Build-RunsXml($Segments) {\n    $runs = \"\"\n    foreach ($seg in $Segments) {\n        $text = $seg[0]\n        $flagged = $seg[1]\n        $run = \"<w:r><w:rPr><w:lang w:val=`\"en-US`\"/></w:rPr><w:t>$text</w:t></w:r>\"\n        if ($flagged) {\n            $runs += \"<w:proofErr w:type=`\"spellStart`\"/>$run<w:proofErr w:type=`\"spellEnd`\"/>\"\n        } else {\n            $runs += $run\n        }\n    }\n    return $runs\n}\n\nfunction Build-PackageXml($RunsXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' + $RunsXml + '</w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# Find the (1-based) index of the paragraph immediately AFTER the one whose\n# text starts with $AnchorText. Done up front, before any paragraph content\n# is modified, and matched case-sensitively so it can never accidentally\n# match text this script itself is about to insert.\nfunction Find-TargetParagraphIndex($Doc, $AnchorText) {\n    $paras = $Doc.Paragraphs\n    $count = $paras.Count\n    for ($i = 1; $i -le $count; $i++) {\n        if ($paras.Item($i).Range.Text -clike ($AnchorText + \"*\")) {\n            return $i + 1\n        }\n    }\n    return -1\n}\n\n# Insert $Segments' runs at the very start of paragraph number $ParaIndex\n# (1-based), leaving that paragraph's own mark / pPr untouched. A\n# zero-length range at the paragraph's Start is used (rather than the\n# paragraph's full Range) so InsertXML merges the new runs into the\n# existing, empty paragraph instead of splitting it into a new one.\nfunction Insert-RunsIntoParagraph($Doc, $ParaIndex, $Segments) {\n    $target = $Doc.Paragraphs.Item($ParaIndex)\n    $insertPoint = $Doc.Range($target.Range.Start, $target.Range.Start)\n    $insertPoint.InsertXML((Build-PackageXml (Build-RunsXml $Segments)))\n}\n\n# Resolve all three target paragraphs before mutating anything.\n$targetDepartments = Find-TargetParagraphIndex $d \"ALTER TABLE DEPARTMENTS ADD FOREIGN KEY\"\n$targetEmployees    = Find-TargetParagraphIndex $d \"ALTER TABLE EMPLOYEES ADD FOREIGN KEY\"\n$targetJobHistory   = Find-TargetParagraphIndex $d \"ALTER TABLE JOB_HISTORY ADD FOREIGN KEY\"\n\nInsert-RunsIntoParagraph $d $targetDepartments @(\n    ,@(\"alter table DEPARTMENTS ADD FOREIGN KEY(\", $false)\n    ,@(\"location_id\", $true)\n    ,@(\") REFERENCES LOCATIONS(\", $false)\n    ,@(\"location_id\", $true)\n    ,@(\");\", $false)\n)\n\nInsert-RunsIntoParagraph $d $targetEmployees @(\n    ,@(\"alter table JOB_HISTORY ADD FOREIGN KEY(\", $false)\n    ,@(\"employee_id\", $true)\n    ,@(\") REFERENCES EMPLOYEES(\", $false)\n    ,@(\"employee_id\", $true)\n    ,@(\");\", $false)\n)\n\nInsert-RunsIntoParagraph $d $targetJobHistory @(\n    ,@(\"alter table JOB_HISTORY ADD FOREIGN KEY(\", $false)\n    ,@(\"job_id\", $true)\n    ,@(\") REFERENCES JOBS(\", $false)\n    ,@(\"job_id\", $true)\n    ,@(\");\", $false)\n)\n"}
